$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so Excel keeps them as text (preserving exact formatting, e.g. trailing zeros)
# instead of auto-converting the assigned string into a floating point number.
$textCells = @("D4","D5","D6","D7","D9","D11","D12","D14","D20","D21","D22","D23","D24","D25","D28","D29","D30","D31","D32","D35","D36","D37","D38","D39","D40","D41","D43","D44","D47","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "66.103.60"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "3.177.40"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "605.25"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").Value = "154.82"
$ws.Range("E6").Value = "  +0.93%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.176.60"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("D9").Value = "0.546"
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").Value = "5.70"
$ws.Range("E11").Value = "  -6.70%  "
$ws.Range("D12").Value = "0.517"
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("D14").Value = "38.28"
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("D15").Value = "3.699.04"
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("D16").Value = "66.130.62"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").Value = "3.178.98"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("D20").Value = "508.74"
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("D21").Value = "15.36"
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("D22").Value = "0.729"
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("D23").Value = "8.01"
$ws.Range("E23").Value = "  -1.69%  "
$ws.Range("D24").Value = "14.79"
$ws.Range("E24").Value = "  -3.78%  "
$ws.Range("D25").Value = "84.44"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("D28").Value = "9.16"
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("D29").Value = "2.38"
$ws.Range("E29").Value = "  +4.64%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "7.19"
$ws.Range("E30").Value = "  +4.66%  "
$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").Value = "2.99"
$ws.Range("E31").Value = "  +4.12%  "
$ws.Range("D32").Value = "27.95"
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  -3.23%  "
$ws.Range("D35").Value = "6.50"
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("D36").Value = "505.13"
$ws.Range("E36").Value = "  +3.83%  "
$ws.Range("D37").Value = "55.17"
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("D38").Value = "0.0878"
$ws.Range("E38").Value = "  -2.96%  "
$ws.Range("D39").Value = "0.0419"
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("D40").Value = "0.128"
$ws.Range("E40").Value = "  +5.84%  "
$ws.Range("D41").Value = "8.78"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("D42").Value = "0.0₃0685"
$ws.Range("E42").Value = "  +5.85%  "
$ws.Range("D43").Value = "2.85"
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("D44").Value = "0.297"
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").Value = "2.831.06"
$ws.Range("E46").Value = "  -4.14%  "
$ws.Range("D47").Value = "28.05"
$ws.Range("E47").Value = "  -1.78%  "
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "2.36"
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").Value = "2.61"
$ws.Range("E51").Value = "  +6.25%  "
